$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.673.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.740.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -10.10%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4896'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -7.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.47'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -7.96%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -18.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06041'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -12.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.743.51'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06832'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -12.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -20.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.421'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -12.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.85'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -15.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5661'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -25.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '25.709.62'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.27'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -19.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006560'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -17.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.963.12'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.025'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -13.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.887'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -15.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.019'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -16.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.58'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.468'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -13.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.807'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -17.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.66'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -13.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '101.57'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -8.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.714'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -13.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07943'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -9.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.369'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -17.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04377'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -9.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.0000'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.621'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -11.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9678'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -14.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5907'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -19.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.662'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -14.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01505'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -12.83%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.99'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.92%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.869'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -19.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.142'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -12.82%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3715'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -22.73%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7237'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -19.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05216'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -10.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1070'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -13.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.93'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -14.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '51.78'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -14.19%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9998'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.24%  '
